$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: replace the "Example 2" placeholder with the real log entry
$ws.Range("A5").Value = 'Working on the requirements'
$ws.Range("B5").Value = 4.25
$ws.Range("C5").Value = [datetime]"2022-06-18"
$ws.Range("D5").Value = 'After working on the first route to fetch all the facilities, my understanding of rest API''s had greatly improved. At first, I forgot to implement the tags of the facility. When the user called /facility, they got to see all the facilities but without their corresponding tags.  After this step, it was pretty easy to replicate the other CRUD routes. The POST was important because it had to be easy to use. The user should be able to create a facility with its tags and location. If the tags or location do not exist, they will be created. The PUT is also very important because this has to update the right facility along with it''s tags. The DELETE route was pretty simple and straight forward. After getting the basics working, it was time to start on the filter method. I gave this a lot of thought but I think that the way I prepare my statements to build the query depending on the filter criteria is safe against SQL injections. Every input the user can give will be sanitized. '
$ws.Range("D5").WrapText = $true

# Row 6: replace the "Example 3" placeholder with the real log entry, clear the bonus "x" marker
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = [datetime]"2022-06-26"
$ws.Range("D6").Value = '
Last week was not very productive for my assignment. During the week I work at my part time job. This is the company where I had my internship at the time of my MBO course. It was my intention to continue with the assignment on the Friday after work but instead I got sick. It was not Corona but a flu of some kind. Luckily for me, I have improved a lot in the last few days. During the last 2 days I have made some improvements to my application. I reflected on the way I handle errors. This could be done without the try catch method because this was a very cheap answer. To improve this I anticipated what parts of the code could give an error and tried to minimize the result. After the optimizations I am planning on implementing the bonus features. I have to keep in mind that it has already been a week since I’ve started on the assignment and that I should not take too much time. '
$ws.Range("D6").WrapText = $true
$ws.Range("A6").Value = 'Improving my code'
$ws.Range("E6").ClearContents()

# Leave the selection on the next empty row, like the author did before saving
$ws.Range("D7").Select()
